# Datum_Conversion.xlsx update — CPRA postproc re-run pulling in additional
# USACE/USGS gauge stations.
#
# Changes applied:
#   1. Row 7 (station that used to be "82740") is replaced by station
#      "82742" and its offset reset to 0 (new/unreviewed station).
#   2. Three new station rows are appended at the bottom of the table
#      (rows 31-33): "01480", "76560" and "073814675".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Station in row 7 changed from 82740 -> 82742, offset reset to 0 ---
$ws.Cells.Item(7, 2).NumberFormat = "@"
$ws.Cells.Item(7, 2).Value = "82742"
$ws.Cells.Item(7, 3).NumberFormat = "0.00"
$ws.Cells.Item(7, 3).Value = 0

# --- 2. New station rows appended after the existing last row (30) ---

# Row 31: USACE station 01480, offset 0
$ws.Cells.Item(31, 1).Value = "USACE"
$ws.Cells.Item(31, 2).NumberFormat = "@"
$ws.Cells.Item(31, 2).Value = "01480"
$ws.Cells.Item(31, 3).NumberFormat = "0.00"
$ws.Cells.Item(31, 3).Value = 0

# Row 32: USACE station 76560, offset -1.214
$ws.Cells.Item(32, 1).Value = "USACE"
$ws.Cells.Item(32, 2).Value = 76560
$ws.Cells.Item(32, 2).NumberFormat = "@"
$ws.Cells.Item(32, 3).NumberFormat = "0.00"
$ws.Cells.Item(32, 3).Value = -1.214

# Row 33: USGS station 073814675, offset 0
$ws.Cells.Item(33, 1).Value = "USGS"
$ws.Cells.Item(33, 2).NumberFormat = "@"
$ws.Cells.Item(33, 2).Value = "073814675"
$ws.Cells.Item(33, 3).NumberFormat = "0.00"
$ws.Cells.Item(33, 3).Value = 0

# Reflect the new selection/dimension as left by the author's last save.
$ws.Range("C8").Select() | Out-Null
